# Applies the "added PV capacity and temperature factors" commit.
$wb = $excel.ActiveWorkbook

$wsSets = $wb.Worksheets.Item(1)          # "Sets"
$wsGeneral = $wb.Worksheets.Item(2)       # "General Data"
$wsNewInv = $wb.Worksheets.Item(3)        # "Costs new investments"
$wsDefaultSys = $wb.Worksheets.Item(4)    # "Costs default system"
$wsIrrTemp = $wb.Worksheets.Item(6)       # "Irradiation and temperatur"

# ---------------------------------------------------------------------
# Sheet "General Data": new parameter rows 19-24
# ---------------------------------------------------------------------
$wsGeneral.Range("A19").Value = "Bonus shifting"
$wsGeneral.Range("B19").Value = "p_shifting"
$wsGeneral.Range("C19").Value = 0.03

$wsGeneral.Range("A20").Value = "Irradiation STC"
$wsGeneral.Range("B20").Value = "Irr_STC [kW/m²]"
$wsGeneral.Range("C20").Value = 1

$wsGeneral.Range("A21").Value = "Temperature STC"
$wsGeneral.Range("B21").Value = "T_STC [°C]"
$wsGeneral.Range("C21").Value = 25

$wsGeneral.Range("A22").Value = "Temperatur factor PV"
$wsGeneral.Range("B22").Value = "[-%/°C]"
$wsGeneral.Range("C22").Value = 0.3

$wsGeneral.Range("A23").Value = "Performance ratio PV"
$wsGeneral.Range("B23").Value = "PR"
$wsGeneral.Range("C23").Value = 0.85

$wsGeneral.Range("A24").Value = "Surface Factor PV"
$wsGeneral.Range("B24").Value = "30 grad, süd"
$wsGeneral.Range("C24").Value = 1.1

# ---------------------------------------------------------------------
# Sheet "Sets": new PV/ST/Grid/Car/Battery/HP flow-mapping columns H:M
# ---------------------------------------------------------------------
$wsSets.Range("H1").Value = "PV to"
$wsSets.Range("H2").Value = "Car"
$wsSets.Range("H3").Value = "Grid"
$wsSets.Range("H4").Value = "Battery"
$wsSets.Range("H5").Value = "Household"
$wsSets.Range("H6").Value = "HP"
$wsSets.Range("H7").Value = "Curtailment"

$wsSets.Range("I1").Value = "ST to"
$wsSets.Range("I2").Value = "DH"
$wsSets.Range("I3").Value = "Household"

$wsSets.Range("J1").Value = "Electric Grid to"
$wsSets.Range("J2").Value = "HP"
$wsSets.Range("J3").Value = "Household"
$wsSets.Range("J4").Value = "Car"
$wsSets.Range("J5").Value = "Battery"

$wsSets.Range("K1").Value = "Car to"
$wsSets.Range("K2").Value = "Battery"

$wsSets.Range("L1").Value = "Battery to"
$wsSets.Range("L2").Value = "Car"

$wsSets.Range("M1").Value = "HP to"
$wsSets.Range("M2").Value = "Household"

# ---------------------------------------------------------------------
# Cosmetic view-state updates (selection, scroll, active sheet/tab)
# ---------------------------------------------------------------------

# "Costs default system": column E bestfit-ish width + new selection
$wsDefaultSys.Columns.Item(5).ColumnWidth = 14.7265625
$wsDefaultSys.Range("E1").Select()

# "Costs new investments": selection moves, loses tab-selected flag
$wsNewInv.Range("D18").Select()

# "Irradiation and temperatur": selection changes
$wsIrrTemp.Range("G16").Select()

# "Sets": selection changes
$wsSets.Range("I1").Select()

# "General Data": becomes the active tab, with a new selection/top row
$wsGeneral.Activate()
$excel.ActiveWindow.ScrollRow = 4
$wsGeneral.Range("F22").Select()
